{"js": "// Update L/FNG Review Form Labels\n// Replace the \"Resolution for Application to Proceed to the ALC\" table\n// label with \"What is the outcome of the Board/Council resolution?\"\n\nconst body = context.document.body;\nconst results = body.search(\"Resolution for Application to Proceed to the ALC\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"What is the outcome of the Board/Council resolution?\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update L/FNG Review Form Labels\n# Replace the \"Resolution for Application to Proceed to the ALC\" table\n# label with \"What is the outcome of the Board/Council resolution?\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Resolution for Application to Proceed to the ALC\"\n$find.Replacement.Text = \"What is the outcome of the Board/Council resolution?\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n"}
